$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 8
$ws.Range("G4").Value = 2.6
$ws.Range("I4").Value = 2.88
$ws.Range("J4").Value = 1.08
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 1.4
$ws.Range("M4").Value = 2.75
$ws.Range("N4").Value = 2.3
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.91
$ws.Range("S4").Value = 1.8
$ws.Range("T4").Value = 7.5
$ws.Range("U4").Value = 12
$ws.Range("V4").Value = 10
$ws.Range("W4").Value = 26
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 34
$ws.Range("Z4").Value = 8
$ws.Range("AA4").Value = 6
$ws.Range("AB4").Value = 17
$ws.Range("AC4").Value = 51
$ws.Range("AD4").Value = 351
$ws.Range("AE4").Value = 8
$ws.Range("AF4").Value = 13
$ws.Range("AG4").Value = 11
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 26
$ws.Range("AJ4").Value = 41
$ws.Range("G5").Value = 3.7
$ws.Range("I5").Value = 1.85
$ws.Range("N5").Value = 1.57
$ws.Range("O5").Value = 2.35
$ws.Range("U5").Value = 21
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 8
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 34
$ws.Range("AD5").Value = 101
$ws.Range("G6").Value = 4.33
$ws.Range("H6").Value = 3.9
$ws.Range("I6").Value = 1.73
$ws.Range("V6").Value = 15
$ws.Range("Y6").Value = 34
$ws.Range("G8").Value = 2.45
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 1.11
$ws.Range("K8").Value = 6.5
$ws.Range("V8").Value = 11
$ws.Range("AA8").Value = 5.5
$ws.Range("AJ8").Value = 51
$ws.Range("G9").Value = 2.4
$ws.Range("I9").Value = 3.3
$ws.Range("J9").Value = 1.11
$ws.Range("K9").Value = 6.5
$ws.Range("P9").Value = 1.62
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 2.2
$ws.Range("S9").Value = 1.62
$ws.Range("V9").Value = 11
$ws.Range("AA9").Value = 6
$ws.Range("AE9").Value = 7.5
$ws.Range("G10").Value = 1.95
$ws.Range("I10").Value = 4.33
$ws.Range("U10").Value = 8
$ws.Range("AF10").Value = 21
$ws.Range("G11").Value = 2.05
$ws.Range("I11").Value = 3.5
$ws.Range("L11").Value = 1.36
$ws.Range("M11").Value = 3
$ws.Range("N11").Value = 2.2
$ws.Range("O11").Value = 1.65
$ws.Range("R11").Value = 1.91
$ws.Range("S11").Value = 1.8
$ws.Range("W11").Value = 19
$ws.Range("X11").Value = 19
$ws.Range("Z11").Value = 8.5
$ws.Range("AD11").Value = 351
$ws.Range("AE11").Value = 9.5
$ws.Range("AF11").Value = 17
$ws.Range("G12").Value = 1.48
$ws.Range("H12").Value = 3.9
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 1.05
$ws.Range("K12").Value = 11
$ws.Range("X12").Value = 13
$ws.Range("AD12").Value = 301
$ws.Range("AG12").Value = 19
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = 2.3
$ws.Range("J15").Value = 1.06
$ws.Range("K15").Value = 10
$ws.Range("R15").Value = 1.73
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 9.5
$ws.Range("V15").Value = 11
$ws.Range("X15").Value = 23
$ws.Range("AG15").Value = 9.5
$ws.Range("G17").Value = 1.87
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 3.4
$ws.Range("M17").Value = 3.8
$ws.Range("O17").Value = 1.98
$ws.Range("R17").Value = 1.64
$ws.Range("S17").Value = 2.13
$ws.Range("T17").Value = 7.3
$ws.Range("U17").Value = 8.25
$ws.Range("V17").Value = 7.2
$ws.Range("W17").Value = 13.5
$ws.Range("X17").Value = 11.5
$ws.Range("AA17").Value = 6.3
$ws.Range("AB17").Value = 11
$ws.Range("AC17").Value = 40
$ws.Range("AD17").Value = 250
$ws.Range("AE17").Value = 10.25
$ws.Range("AF17").Value = 16
$ws.Range("AG17").Value = 10
$ws.Range("AI17").Value = 22
$ws.Range("AJ17").Value = 25
$ws.Range("G18").Value = 2.38
$ws.Range("I18").Value = 3
